# LOQ4261.docx edit script
# Applies a rotation of text contents across several paragraphs (Objetivos,
# Docente, Programa resumido, Programa, Avaliacao), removes the long
# "Norma de recuperacao" value run, and relocates the "Bibliografia"
# heading (+ the "8971158 - Claudemir Leif Tramarico" paragraph) from
# right after "Avaliacao" to right after the bibliography list, just
# before "Requisitos".

$d = $word.ActiveDocument

# --- 1) Objetivos (PT) value paragraph: was "Apresentar um quadro..."
#        becomes the old "Programa resumido" (PT) text.
$d.Paragraphs(6).Range.Text = "Caracterização do planejamento e controle da produção. 2. Planejamento agregado da produção. 3. Planejamento mestre da produção. 4. Planejamento e controle de estoques. 5. Planejamento de recursos de materiais (MRP). 6. Programação detalhada da produção. 7. Sistema MRPII e Sistema ERP. 8.Tambor-Pulmão-Corda - OPT. 9. Teoria das Restrições (TOC)."

# --- 2) Objetivos (EN, italic) value paragraph: was "To present a conceptual..."
#        becomes the old "Programa resumido" (EN) text.
$d.Paragraphs(7).Range.Text = "Characterization of production programming and control. 2. Aggregate Production Planning. 3. Master Production Schedulling. 4. Inventory planning and control. 5. Material Requirement Planning (MRP). 6. Detailed scheduling of production. 7. Production control systems. 8. Drum-Buffer-Rope – Opt; 9. Theory of Constraints (TOC)"

# --- 3) Docente(s) value paragraph (ListBullet): was "8971158 - Claudemir..."
#        becomes the OLD "Objetivos" (PT) text.
$d.Paragraphs(9).Range.Text = "Apresentar um quadro conceitual de análise para auxiliar na formulação, avaliação e desenvolvimento de modelos para Planejamento, Programação e Controle da Produção nos diferentes ambientes de produção."

# --- 4) "Programa resumido" (EN, italic) value paragraph: was "Characterization of production..."
#        becomes the OLD "Objetivos" (EN) text.
$d.Paragraphs(12).Range.Text = "To present a conceptual framework of analysis to assist in the formulation, evaluation and development of models for Planning, Programming and Production Control in different production environments."

# --- 5) "Programa" (PT) value paragraph: was "Caracterização do planejamento..."
#        becomes the OLD "Avaliacao" Metodo value text.
$d.Paragraphs(14).Range.Text = "Provas, atividades em grupo e atividades individuais."

# (Paragraph 11 "Programa resumido" (PT) and paragraph 15 "Programa" (EN)
#  keep their original text unchanged.)

# --- 6) Avaliacao block (single ListBullet paragraph with several runs
#        separated by <w:br/>). Update the "Metodo" and "Criterio" values,
#        scoped to that paragraph to avoid touching the (currently)
#        duplicate text elsewhere in the document.
$pAval = $d.Paragraphs(17).Range.Duplicate
$pAval.Find.Execute("Provas, atividades em grupo e atividades individuais.", $true, $false, $false, $false, $false, $true, 1, $false, "Média das atividades avaliativas", 2) | Out-Null

$pAval2 = $d.Paragraphs(17).Range.Duplicate
$pAval2.Find.Execute("Média das atividades avaliativas", $true, $false, $false, $false, $false, $true, 1, $false, "MF = (0,5 M + 0,5 R) M = Média de aproveitamento do aluno, antes da recuperação R = Nota de uma prova de recuperação MF = nota final de aproveitamento, após a recuperação Aprovação com média final de aproveitamento maior ou igual a 5,0. A recuperação deverá consistir em uma prova escrita englobando a matéria toda do semestre. Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.", 2) | Out-Null

# --- 7) Remove the (now stale) "Norma de recuperacao" value run entirely
#        (its whole text is matched so the whole run is consumed/deleted).
$pAval3 = $d.Paragraphs(17).Range.Duplicate
$pAval3.Find.Execute("MF = (0,5 M + 0,5 R) M = Média de aproveitamento do aluno, antes da recuperação R = Nota de uma prova de recuperação MF = nota final de aproveitamento, após a recuperação Aprovação com média final de aproveitamento maior ou igual a 5,0. A recuperação deverá consistir em uma prova escrita englobando a matéria toda do semestre. Terá direito à prova de recuperação aqueles alunos reprovados com nota acima de 3,0 e frequência mínima de 70%.", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# --- 8) Remove the "Bibliografia" Heading2 paragraph that currently sits
#        right after the Avaliacao block (paragraph 18), then merge the
#        Avaliacao paragraph with the bibliography list paragraph that
#        follows it (deleting the paragraph mark between them), restoring
#        the ListBullet style on the resulting merged paragraph.
$pBibHeading = $d.Paragraphs(18)
$pBibHeading.Range.Delete()

$pAvalEnd = $d.Paragraphs(17)
$mark = $d.Range($pAvalEnd.Range.End - 1, $pAvalEnd.Range.End)
$mark.Delete()

$merged = $d.Paragraphs(17)
$merged.Style = "ListBullet"

# --- 9) Re-insert the "Bibliografia" heading and the "8971158 - Claudemir
#        Leif Tramarico" paragraph right after the bibliography list
#        (i.e. right before "Requisitos").
$pReq = $d.Paragraphs(18)
$insertPoint = $d.Range($pReq.Range.Start, $pReq.Range.Start)
$insertPoint.InsertParagraphBefore()
$insertPoint.InsertParagraphBefore()

$pNew1 = $d.Paragraphs(18)
$pNew1.Range.Text = "Bibliografia"
$pNew1.Style = "Heading2"

$pNew2 = $d.Paragraphs(19)
$pNew2.Range.Text = "8971158 - Claudemir Leif Tramarico"

Write-Host "Edit complete"
